$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 165 (last existing data row) into new rows 166-174
$ws.Range("A165:AC165").Copy()
$ws.Range("A166:AC174").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 166
$ws.Cells.Item(166,1).Value2 = 164
$ws.Cells.Item(166,2).Value2 = 6818360
$ws.Cells.Item(166,3).Value2 = "Hungary NB I"
$ws.Cells.Item(166,4).Value2 = "Hungary NB I"
$ws.Cells.Item(166,5).Value2 = 45395.60416666666
$ws.Cells.Item(166,6).Value2 = "Paksi"
$ws.Cells.Item(166,7).Value2 = "MOL Fehervar FC"
$ws.Cells.Item(166,8).Value2 = 1
$ws.Cells.Item(166,9).Value2 = 2
$ws.Cells.Item(166,10).Value2 = "A"
$ws.Cells.Item(166,11).Value2 = 1.666
$ws.Cells.Item(166,12).Value2 = 3.5
$ws.Cells.Item(166,13).Value2 = 5.5
$ws.Cells.Item(166,14).Value2 = 1.909
$ws.Cells.Item(166,15).Value2 = 3.4
$ws.Cells.Item(166,16).Value2 = 4
$ws.Cells.Item(166,17).Value2 = -0.5
$ws.Cells.Item(166,18).Value2 = 1.875
$ws.Cells.Item(166,19).Value2 = 1.975
$ws.Cells.Item(166,20).Value2 = 2.75
$ws.Cells.Item(166,21).Value2 = 1.825
$ws.Cells.Item(166,22).Value2 = 2.025
$ws.Cells.Item(166,23).Value2 = -1
$ws.Cells.Item(166,24).Value2 = -1
$ws.Cells.Item(166,25).Value2 = 3
$ws.Cells.Item(166,26).Value2 = -1
$ws.Cells.Item(166,27).Value2 = 0.9750000000000001
$ws.Cells.Item(166,28).Value2 = 0.4125
$ws.Cells.Item(166,29).Value2 = -0.5

# Row 167
$ws.Cells.Item(167,1).Value2 = 165
$ws.Cells.Item(167,2).Value2 = 6820520
$ws.Cells.Item(167,3).Value2 = "Hungary NB I"
$ws.Cells.Item(167,4).Value2 = "Hungary NB I"
$ws.Cells.Item(167,5).Value2 = 45396.375
$ws.Cells.Item(167,6).Value2 = "MTK Budapest"
$ws.Cells.Item(167,7).Value2 = "Diosgyori VTK"
$ws.Cells.Item(167,8).Value2 = 1
$ws.Cells.Item(167,9).Value2 = 1
$ws.Cells.Item(167,10).Value2 = "D"
$ws.Cells.Item(167,11).Value2 = 1.727
$ws.Cells.Item(167,12).Value2 = 3.4
$ws.Cells.Item(167,13).Value2 = 4.5
$ws.Cells.Item(167,14).Value2 = 2.2
$ws.Cells.Item(167,15).Value2 = 3.2
$ws.Cells.Item(167,16).Value2 = 2.9
$ws.Cells.Item(167,17).Value2 = -0.25
$ws.Cells.Item(167,18).Value2 = 1.975
$ws.Cells.Item(167,19).Value2 = 1.875
$ws.Cells.Item(167,20).Value2 = 3
$ws.Cells.Item(167,21).Value2 = 1.975
$ws.Cells.Item(167,22).Value2 = 1.875
$ws.Cells.Item(167,23).Value2 = -1
$ws.Cells.Item(167,24).Value2 = 2.2
$ws.Cells.Item(167,25).Value2 = -1
$ws.Cells.Item(167,26).Value2 = -0.5
$ws.Cells.Item(167,27).Value2 = 0.4375
$ws.Cells.Item(167,28).Value2 = -1
$ws.Cells.Item(167,29).Value2 = 0.875

# Row 168
$ws.Cells.Item(168,1).Value2 = 166
$ws.Cells.Item(168,2).Value2 = 6860904
$ws.Cells.Item(168,3).Value2 = "Hungary NB I"
$ws.Cells.Item(168,4).Value2 = "Hungary NB I"
$ws.Cells.Item(168,5).Value2 = 45396.47916666666
$ws.Cells.Item(168,6).Value2 = "Zalaegerszegi TE"
$ws.Cells.Item(168,7).Value2 = "Ferencvarosi TC"
$ws.Cells.Item(168,8).Value2 = 2
$ws.Cells.Item(168,9).Value2 = 3
$ws.Cells.Item(168,10).Value2 = "A"
$ws.Cells.Item(168,11).Value2 = 5.5
$ws.Cells.Item(168,12).Value2 = 4
$ws.Cells.Item(168,13).Value2 = 1.5
$ws.Cells.Item(168,14).Value2 = 6.5
$ws.Cells.Item(168,15).Value2 = 3.8
$ws.Cells.Item(168,16).Value2 = 1.45
$ws.Cells.Item(168,17).Value2 = 1.25
$ws.Cells.Item(168,18).Value2 = 1.8
$ws.Cells.Item(168,19).Value2 = 2.05
$ws.Cells.Item(168,20).Value2 = 2.75
$ws.Cells.Item(168,21).Value2 = 1.875
$ws.Cells.Item(168,22).Value2 = 1.975
$ws.Cells.Item(168,23).Value2 = -1
$ws.Cells.Item(168,24).Value2 = -1
$ws.Cells.Item(168,25).Value2 = 0.45
$ws.Cells.Item(168,26).Value2 = 0.4
$ws.Cells.Item(168,27).Value2 = -0.5
$ws.Cells.Item(168,28).Value2 = 0.875
$ws.Cells.Item(168,29).Value2 = -1

# Row 169
$ws.Cells.Item(169,1).Value2 = 167
$ws.Cells.Item(169,2).Value2 = 6818365
$ws.Cells.Item(169,3).Value2 = "Hungary NB I"
$ws.Cells.Item(169,4).Value2 = "Hungary NB I"
$ws.Cells.Item(169,5).Value2 = 45402.39583333334
$ws.Cells.Item(169,6).Value2 = "Puskas Academy"
$ws.Cells.Item(169,7).Value2 = "Paksi"
$ws.Cells.Item(169,8).Clear()
$ws.Cells.Item(169,9).Clear()
$ws.Cells.Item(169,10).Clear()
$ws.Cells.Item(169,11).Value2 = 2
$ws.Cells.Item(169,12).Value2 = 3.4
$ws.Cells.Item(169,13).Value2 = 3.3
$ws.Cells.Item(169,14).Value2 = 1.909
$ws.Cells.Item(169,15).Value2 = 3.5
$ws.Cells.Item(169,16).Value2 = 3.6
$ws.Cells.Item(169,17).Value2 = -0.5
$ws.Cells.Item(169,18).Value2 = 1.925
$ws.Cells.Item(169,19).Value2 = 1.925
$ws.Cells.Item(169,20).Value2 = 2.75
$ws.Cells.Item(169,21).Value2 = 2
$ws.Cells.Item(169,22).Value2 = 1.85
$ws.Cells.Item(169,23).Value2 = 0
$ws.Cells.Item(169,24).Value2 = 0
$ws.Cells.Item(169,25).Value2 = 0
$ws.Cells.Item(169,26).Value2 = 0
$ws.Cells.Item(169,27).Value2 = 0
$ws.Cells.Item(169,28).Clear()
$ws.Cells.Item(169,29).Clear()

# Row 170
$ws.Cells.Item(170,1).Value2 = 168
$ws.Cells.Item(170,2).Value2 = 6818362
$ws.Cells.Item(170,3).Value2 = "Hungary NB I"
$ws.Cells.Item(170,4).Value2 = "Hungary NB I"
$ws.Cells.Item(170,5).Value2 = 45402.5
$ws.Cells.Item(170,6).Value2 = "Ferencvarosi TC"
$ws.Cells.Item(170,7).Value2 = "Kisvarda FC"
$ws.Cells.Item(170,8).Clear()
$ws.Cells.Item(170,9).Clear()
$ws.Cells.Item(170,10).Clear()
$ws.Cells.Item(170,11).Value2 = 1.25
$ws.Cells.Item(170,12).Value2 = 5.5
$ws.Cells.Item(170,13).Value2 = 9
$ws.Cells.Item(170,14).Value2 = 1.181
$ws.Cells.Item(170,15).Value2 = 5.75
$ws.Cells.Item(170,16).Value2 = 13
$ws.Cells.Item(170,17).Value2 = -2
$ws.Cells.Item(170,18).Value2 = 2.05
$ws.Cells.Item(170,19).Value2 = 1.8
$ws.Cells.Item(170,20).Value2 = 3
$ws.Cells.Item(170,21).Value2 = 1.875
$ws.Cells.Item(170,22).Value2 = 1.975
$ws.Cells.Item(170,23).Value2 = 0
$ws.Cells.Item(170,24).Value2 = 0
$ws.Cells.Item(170,25).Value2 = 0
$ws.Cells.Item(170,26).Value2 = 0
$ws.Cells.Item(170,27).Value2 = 0
$ws.Cells.Item(170,28).Clear()
$ws.Cells.Item(170,29).Clear()

# Row 171
$ws.Cells.Item(171,1).Value2 = 169
$ws.Cells.Item(171,2).Value2 = 6818364
$ws.Cells.Item(171,3).Value2 = "Hungary NB I"
$ws.Cells.Item(171,4).Value2 = "Hungary NB I"
$ws.Cells.Item(171,5).Value2 = 45402.60416666666
$ws.Cells.Item(171,6).Value2 = "Diosgyori VTK"
$ws.Cells.Item(171,7).Value2 = "Debreceni VSC"
$ws.Cells.Item(171,8).Clear()
$ws.Cells.Item(171,9).Clear()
$ws.Cells.Item(171,10).Clear()
$ws.Cells.Item(171,11).Value2 = 2.5
$ws.Cells.Item(171,12).Value2 = 3.2
$ws.Cells.Item(171,13).Value2 = 2.625
$ws.Cells.Item(171,14).Value2 = 2.625
$ws.Cells.Item(171,15).Value2 = 3.2
$ws.Cells.Item(171,16).Value2 = 2.5
$ws.Cells.Item(171,17).Value2 = 0
$ws.Cells.Item(171,18).Value2 = 2.025
$ws.Cells.Item(171,19).Value2 = 1.825
$ws.Cells.Item(171,20).Value2 = 2.75
$ws.Cells.Item(171,21).Value2 = 2
$ws.Cells.Item(171,22).Value2 = 1.85
$ws.Cells.Item(171,23).Value2 = 0
$ws.Cells.Item(171,24).Value2 = 0
$ws.Cells.Item(171,25).Value2 = 0
$ws.Cells.Item(171,26).Value2 = 0
$ws.Cells.Item(171,27).Value2 = 0
$ws.Cells.Item(171,28).Clear()
$ws.Cells.Item(171,29).Clear()

# Row 172
$ws.Cells.Item(172,1).Value2 = 170
$ws.Cells.Item(172,2).Value2 = 6818367
$ws.Cells.Item(172,3).Value2 = "Hungary NB I"
$ws.Cells.Item(172,4).Value2 = "Hungary NB I"
$ws.Cells.Item(172,5).Value2 = 45403.37847222222
$ws.Cells.Item(172,6).Value2 = "Kecskemeti TE"
$ws.Cells.Item(172,7).Value2 = "Mezokovesd Zsory"
$ws.Cells.Item(172,8).Clear()
$ws.Cells.Item(172,9).Clear()
$ws.Cells.Item(172,10).Clear()
$ws.Cells.Item(172,11).Value2 = 1.727
$ws.Cells.Item(172,12).Value2 = 3.5
$ws.Cells.Item(172,13).Value2 = 4.333
$ws.Cells.Item(172,14).Value2 = 1.666
$ws.Cells.Item(172,15).Value2 = 3.6
$ws.Cells.Item(172,16).Value2 = 4.5
$ws.Cells.Item(172,17).Value2 = -0.75
$ws.Cells.Item(172,18).Value2 = 1.9
$ws.Cells.Item(172,19).Value2 = 1.95
$ws.Cells.Item(172,20).Value2 = 2.5
$ws.Cells.Item(172,21).Value2 = 2.025
$ws.Cells.Item(172,22).Value2 = 1.825
$ws.Cells.Item(172,23).Value2 = 0
$ws.Cells.Item(172,24).Value2 = 0
$ws.Cells.Item(172,25).Value2 = 0
$ws.Cells.Item(172,26).Value2 = 0
$ws.Cells.Item(172,27).Value2 = 0
$ws.Cells.Item(172,28).Clear()
$ws.Cells.Item(172,29).Clear()

# Row 173
$ws.Cells.Item(173,1).Value2 = 171
$ws.Cells.Item(173,2).Value2 = 6818366
$ws.Cells.Item(173,3).Value2 = "Hungary NB I"
$ws.Cells.Item(173,4).Value2 = "Hungary NB I"
$ws.Cells.Item(173,5).Value2 = 45403.45833333334
$ws.Cells.Item(173,6).Value2 = "MOL Fehervar FC"
$ws.Cells.Item(173,7).Value2 = "Zalaegerszegi TE"
$ws.Cells.Item(173,8).Clear()
$ws.Cells.Item(173,9).Clear()
$ws.Cells.Item(173,10).Clear()
$ws.Cells.Item(173,11).Value2 = 1.909
$ws.Cells.Item(173,12).Value2 = 3.4
$ws.Cells.Item(173,13).Value2 = 3.6
$ws.Cells.Item(173,14).Value2 = 1.909
$ws.Cells.Item(173,15).Value2 = 3.4
$ws.Cells.Item(173,16).Value2 = 3.6
$ws.Cells.Item(173,17).Value2 = -0.5
$ws.Cells.Item(173,18).Value2 = 1.975
$ws.Cells.Item(173,19).Value2 = 1.875
$ws.Cells.Item(173,20).Value2 = 2.75
$ws.Cells.Item(173,21).Value2 = 2
$ws.Cells.Item(173,22).Value2 = 1.85
$ws.Cells.Item(173,23).Value2 = 0
$ws.Cells.Item(173,24).Value2 = 0
$ws.Cells.Item(173,25).Value2 = 0
$ws.Cells.Item(173,26).Value2 = 0
$ws.Cells.Item(173,27).Value2 = 0
$ws.Cells.Item(173,28).Clear()
$ws.Cells.Item(173,29).Clear()

# Row 174
$ws.Cells.Item(174,1).Value2 = 172
$ws.Cells.Item(174,2).Value2 = 6818363
$ws.Cells.Item(174,3).Value2 = "Hungary NB I"
$ws.Cells.Item(174,4).Value2 = "Hungary NB I"
$ws.Cells.Item(174,5).Value2 = 45403.5625
$ws.Cells.Item(174,6).Value2 = "Ujpest"
$ws.Cells.Item(174,7).Value2 = "MTK Budapest"
$ws.Cells.Item(174,8).Clear()
$ws.Cells.Item(174,9).Clear()
$ws.Cells.Item(174,10).Clear()
$ws.Cells.Item(174,11).Value2 = 2
$ws.Cells.Item(174,12).Value2 = 3.4
$ws.Cells.Item(174,13).Value2 = 3.3
$ws.Cells.Item(174,14).Value2 = 2.1
$ws.Cells.Item(174,15).Value2 = 3.4
$ws.Cells.Item(174,16).Value2 = 3.1
$ws.Cells.Item(174,17).Value2 = -0.25
$ws.Cells.Item(174,18).Value2 = 1.9
$ws.Cells.Item(174,19).Value2 = 1.95
$ws.Cells.Item(174,20).Value2 = 2.75
$ws.Cells.Item(174,21).Value2 = 1.825
$ws.Cells.Item(174,22).Value2 = 2.025
$ws.Cells.Item(174,23).Value2 = 0
$ws.Cells.Item(174,24).Value2 = 0
$ws.Cells.Item(174,25).Value2 = 0
$ws.Cells.Item(174,26).Value2 = 0
$ws.Cells.Item(174,27).Value2 = 0
$ws.Cells.Item(174,28).Clear()
$ws.Cells.Item(174,29).Clear()
